$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Conf1', 'AdaBoost', '0.25 ± 0.11', '1.03 ± 0.07', '61.18 ± 5.13'),
    @('Conf1', 'BaggedDT', '0.24 ± 0.08', '1.04 ± 0.09', '61.65 ± 6.37'),
    @('Conf1', 'BayesianNN', '-0.01 ± 0.09', '1.20 ± 0.10', '71.30 ± 6.76'),
    @('Conf1', 'CatBoost', '0.29 ± 0.07', '1.00 ± 0.07', '59.51 ± 5.05'),
    @('Conf1', 'GBDT', '0.24 ± 0.09', '1.04 ± 0.07', '61.66 ± 4.78'),
    @('Conf1', 'KNN', '-0.02 ± 0.08', '1.20 ± 0.10', '71.54 ± 6.91'),
    @('Conf1', 'LightGBM', '0.22 ± 0.11', '1.05 ± 0.06', '62.29 ± 4.67'),
    @('Conf1', 'MLP', '-0.04 ± 0.04', '1.22 ± 0.10', '72.35 ± 7.31'),
    @('Conf1', 'RF', '0.29 ± 0.04', '1.00 ± 0.07', '59.60 ± 5.36'),
    @('Conf1', 'SVR', '-0.01 ± 0.09', '1.20 ± 0.09', '71.16 ± 6.43'),
    @('Conf1', 'StackEns', '0.25 ± 0.06', '1.03 ± 0.09', '61.40 ± 6.42'),
    @('Conf1', 'XGBoost', '0.23 ± 0.07', '1.04 ± 0.08', '61.96 ± 5.67'),
    @('Conf1', '▸ Average', '0.16 ± 0.08', '1.09 ± 0.08', '64.63 ± 5.90'),
    @('Conf2', 'AdaBoost', '0.24 ± 0.07', '1.03 ± 0.08', '61.53 ± 5.71'),
    @('Conf2', 'BaggedDT', '0.26 ± 0.06', '1.03 ± 0.07', '61.08 ± 5.40'),
    @('Conf2', 'BayesianNN', '0.04 ± 0.10', '1.16 ± 0.08', '69.27 ± 5.29'),
    @('Conf2', 'CatBoost', '0.30 ± 0.05', '1.00 ± 0.07', '59.31 ± 5.23'),
    @('Conf2', 'GBDT', '0.17 ± 0.14', '1.08 ± 0.08', '64.31 ± 5.57'),
    @('Conf2', 'KNN', '0.05 ± 0.11', '1.16 ± 0.07', '68.92 ± 4.45'),
    @('Conf2', 'LightGBM', '0.27 ± 0.07', '1.01 ± 0.06', '60.40 ± 4.83'),
    @('Conf2', 'MLP', '-0.39 ± 0.45', '1.39 ± 0.22', '82.69 ± 13.88'),
    @('Conf2', 'RF', '0.29 ± 0.04', '1.00 ± 0.07', '59.70 ± 5.07'),
    @('Conf2', 'SVR', '0.01 ± 0.11', '1.18 ± 0.07', '70.26 ± 4.90'),
    @('Conf2', 'StackEns', '0.26 ± 0.05', '1.02 ± 0.08', '60.87 ± 6.13'),
    @('Conf2', 'XGBoost', '0.25 ± 0.10', '1.03 ± 0.08', '61.36 ± 5.81'),
    @('Conf2', '▸ Average', '0.15 ± 0.11', '1.09 ± 0.08', '64.97 ± 6.02'),
    @('Conf3', 'AdaBoost', '0.41 ± 0.17', '0.90 ± 0.12', '53.72 ± 6.91'),
    @('Conf3', 'BaggedDT', '0.38 ± 0.19', '0.93 ± 0.12', '55.33 ± 6.86'),
    @('Conf3', 'BayesianNN', '0.35 ± 0.20', '0.95 ± 0.11', '56.42 ± 6.05'),
    @('Conf3', 'CatBoost', '0.42 ± 0.12', '0.90 ± 0.09', '53.53 ± 5.05'),
    @('Conf3', 'GBDT', '0.41 ± 0.15', '0.91 ± 0.08', '53.79 ± 4.12'),
    @('Conf3', 'KNN', '0.01 ± 0.11', '1.18 ± 0.10', '70.26 ± 6.36'),
    @('Conf3', 'LightGBM', '0.39 ± 0.18', '0.92 ± 0.09', '54.58 ± 4.81'),
    @('Conf3', 'MLP', '-0.03 ± 0.01', '1.21 ± 0.10', '71.90 ± 7.52'),
    @('Conf3', 'RF', '0.42 ± 0.13', '0.90 ± 0.09', '53.68 ± 5.52'),
    @('Conf3', 'SVR', '0.13 ± 0.13', '1.11 ± 0.11', '65.82 ± 7.22'),
    @('Conf3', 'StackEns', '0.38 ± 0.13', '0.94 ± 0.11', '55.64 ± 6.64'),
    @('Conf3', 'XGBoost', '0.36 ± 0.18', '0.94 ± 0.09', '55.71 ± 4.80'),
    @('Conf3', '▸ Average', '0.30 ± 0.14', '0.98 ± 0.10', '58.36 ± 5.99'),
    @('Conf4', 'AdaBoost', '0.41 ± 0.16', '0.91 ± 0.10', '53.99 ± 5.55'),
    @('Conf4', 'BaggedDT', '0.38 ± 0.18', '0.93 ± 0.11', '55.20 ± 6.58'),
    @('Conf4', 'BayesianNN', '0.35 ± 0.18', '0.95 ± 0.10', '56.32 ± 5.19'),
    @('Conf4', 'CatBoost', '0.44 ± 0.12', '0.88 ± 0.08', '52.44 ± 4.76'),
    @('Conf4', 'GBDT', '0.41 ± 0.16', '0.91 ± 0.08', '53.87 ± 3.92'),
    @('Conf4', 'KNN', '0.06 ± 0.08', '1.15 ± 0.08', '68.53 ± 5.58'),
    @('Conf4', 'LightGBM', '0.41 ± 0.17', '0.91 ± 0.09', '53.81 ± 4.46'),
    @('Conf4', 'MLP', '-0.03 ± 0.07', '1.21 ± 0.11', '72.00 ± 7.88'),
    @('Conf4', 'RF', '0.42 ± 0.12', '0.91 ± 0.09', '53.84 ± 5.38'),
    @('Conf4', 'SVR', '0.03 ± 0.11', '1.17 ± 0.07', '69.63 ± 4.86'),
    @('Conf4', 'StackEns', '0.38 ± 0.15', '0.93 ± 0.11', '55.53 ± 6.78'),
    @('Conf4', 'XGBoost', '0.36 ± 0.19', '0.94 ± 0.09', '55.88 ± 4.55'),
    @('Conf4', '▸ Average', '0.30 ± 0.14', '0.98 ± 0.09', '58.42 ± 5.46'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
